$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Locator_ids")
$ws2 = $wb.Worksheets.Item("Verification_Entities")

# New locator rows appended to Locator_ids (columns: A = name, B = value)
$locatorRows = @(
    @("Amazon_Departments_Drop_Down_xpath", "xpath=(//span[contains(text(),'Departments')])"),
    @("Amazon_Electronics_Drop_Down_xpath", "xpath=(//span[contains(text(),'Electronics')])"),
    @("Amazon_Headphones_Drop_Down_xpath", "xpath=(//span[contains(text(),'Headphones')])"),
    @("Amazon_Headphones_First_Product_xpath", "xpath=(//a[contains(@class,'a-link-normal')]//img)[1]"),
    @("Amazon_Add_To_Cart_Button_xpath", "xpath=(//input[contains(@id,'add-to-cart-button')])[1]"),
    @("Amazon_Cart_Close_Button_xpath", "xpath=(//a[contains(@class,'close-button')])"),
    @("Amazon_Search_Text_Box_xpath", "xpath=(//input[contains(@id,'searchtextbox')])"),
    @("Amazon_Search_Submit_Button_xpath", "xpath=(//input[contains(@type,'submit')])[1]"),
    @("Amazon_Search_Second_Product_xpath", "xpath=(//a[contains(@class,'a-link-normal')]//img)[2]"),
    @("Amazon_Product_Quantity_Drop_Down_xpath", "xpath=//div[contains(@id,'selectQuantity')]//span[@role]"),
    @("Amazon_Product_Quantity_2_xpath", "xpath=//a[@id='quantity_1']"),
    @("Amazon_Product_Cart_Link_xpath", "xpath=(//a[contains(@id,'nav-cart')])"),
    @("Amazon_Product_Titles_Added_In_Cart_xpath", "xpath=(//span[contains(@class,'sc-product-title')])"),
    @("Amazon_Product_Delete_Buttons_Added_In_Cart_xpath", "xpath=(//input[contains(@name,'submit.delete')])"),
    @("Amazon_Cart_First_Product_Quantity_xpath", "xpath=(//select[contains(@name,'quantity')])[1]"),
    @("Amazon_Cart_Proceed_To_Checkout_Button_xpath", "xpath=(//input[contains(@name,'proceedToCheckout')])"),
    @("Amazon_Checkout_Shipping_Address_xpath", "xpath=(//h1[contains(text(),'shipping address')])")
)

$startRow = 7
for ($i = 0; $i -lt $locatorRows.Count; $i++) {
    $r = $startRow + $i
    $ws1.Range("A$r").Value = $locatorRows[$i][0]
    $ws1.Range("B$r").Value = $locatorRows[$i][1]
}

# Widen column B on Locator_ids sheet (engine rounds ColumnWidth to the
# nearest 1/256 character cell; 73.17 round-trips to a stored width of 74)
$ws1.Columns.Item(2).ColumnWidth = 73.17

# New verification row on Verification_Entities
$ws2.Range("A3").Value = "Verify_Amazon_Added_To_Cart_Text_Message_xpath"
$ws2.Range("B3").Value = "xpath=(//h4[contains(text(),'Added to Cart')])[2]"

# Update selections / active sheet to match final state
$ws1.Range("B11").Select()
$ws2.Range("A3").Select()
$ws2.Activate()
